function Set-CellText($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $style = $cell.Style
    $cell.Value = $val
    $cell.Style = $style
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "'67.536.75"
Set-CellText $ws "E2" "'  +0.99%  "
Set-CellText $ws "D3" "'3.872.03"
Set-CellText $ws "E3" "'  +0.14%  "
Set-CellText $ws "D4" "'1.00"
Set-CellText $ws "E4" "'  +0.08%  "
Set-CellText $ws "D5" "'469.06"
Set-CellText $ws "E5" "'  +10.31%  "
Set-CellText $ws "D6" "'148.46"
Set-CellText $ws "E6" "'  +13.19%  "
Set-CellText $ws "D8" "'0.999"
Set-CellText $ws "E8" "'  +0.05%  "
Set-CellText $ws "D9" "'0.751"
Set-CellText $ws "E9" "'  +3.41%  "
Set-CellText $ws "E10" "'  -2.97%  "
Set-CellText $ws "E11" "'  -9.43%  "
Set-CellText $ws "D12" "'43.85"
Set-CellText $ws "E12" "'  +7.20%  "
Set-CellText $ws "D13" "'10.44"
Set-CellText $ws "E13" "'  +1.86%  "
Set-CellText $ws "D14" "'4.489.80"
Set-CellText $ws "E14" "'  +0.63%  "
Set-CellText $ws "D15" "'14.81"
Set-CellText $ws "E15" "'  -6.64%  "
Set-CellText $ws "D16" "'3.884.00"
Set-CellText $ws "E16" "'  -0.16%  "
Set-CellText $ws "E17" "'  -0.22%  "
Set-CellText $ws "E18" "'  +0.58%  "
Set-CellText $ws "E19" "'  +7.86%  "
Set-CellText $ws "D20" "'67.651.06"
Set-CellText $ws "E20" "'  +1.00%  "
Set-CellText $ws "D21" "'432.60"
Set-CellText $ws "E21" "'  +4.58%  "
Set-CellText $ws "D22" "'14.83"
Set-CellText $ws "E22" "'  -0.45%  "
Set-CellText $ws "D23" "'3.30"
Set-CellText $ws "E23" "'  +8.85%  "
Set-CellText $ws "D24" "'88.72"
Set-CellText $ws "E24" "'  +5.09%  "
Set-CellText $ws "D25" "'3.59"
Set-CellText $ws "E25" "'  +10.78%  "
Set-CellText $ws "D26" "'10.33"
Set-CellText $ws "E26" "'  +13.59%  "
Set-CellText $ws "D27" "'37.74"
Set-CellText $ws "E27" "'  +0.02%  "
Set-CellText $ws "D28" "'10.15"
Set-CellText $ws "E28" "'  +2.02%  "
Set-CellText $ws "D29" "'5.51"
Set-CellText $ws "E29" "'  +3.92%  "
Set-CellText $ws "D30" "'750.71"
Set-CellText $ws "E30" "'  +3.49%  "
Set-CellText $ws "E31" "'  +10.68%  "
Set-CellText $ws "D32" "'13.77"
Set-CellText $ws "E32" "'  +4.66%  "
Set-CellText $ws "D33" "'2.77"
Set-CellText $ws "E33" "'  -0.45%  "
Set-CellText $ws "D34" "'43.25"
Set-CellText $ws "E34" "'  +10.80%  "
Set-CellText $ws "D35" "'0.163"
Set-CellText $ws "E35" "'  +7.26%  "
Set-CellText $ws "D36" "'57.72"
Set-CellText $ws "E36" "'  +3.44%  "
Set-CellText $ws "D37" "'1.00"
Set-CellText $ws "E37" "'  +0.16%  "
Set-CellText $ws "D38" "'5.55"
Set-CellText $ws "E38" "'  +3.11%  "
Set-CellText $ws "D39" "'0.0481"
Set-CellText $ws "E39" "'  +4.23%  "
Set-CellText $ws "D40" "'0.352"
Set-CellText $ws "E40" "'  +12.40%  "
Set-CellText $ws "D41" "'2.91"
Set-CellText $ws "E41" "'  +0.89%  "
Set-CellText $ws "D42" "'2.62"
Set-CellText $ws "E42" "'  +16.58%  "
Set-CellText $ws "E43" "'  +5.65%  "
Set-CellText $ws "E44" "'  -10.37%  "
Set-CellText $ws "E45" "'  +0.10%  "
Set-CellText $ws "E46" "'  +2.22%  "
Set-CellText $ws "B47" "'WEMIXToken"
Set-CellText $ws "C47" "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-CellText $ws "D47" "'2.78"
Set-CellText $ws "E47" "'  +8.15%  "
Set-CellText $ws "B48" "'ApeXProtocol"
Set-CellText $ws "C48" "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-CellText $ws "D48" "'3.24"
Set-CellText $ws "E48" "'  +3.24%  "
Set-CellText $ws "E49" "'  +4.41%  "
Set-CellText $ws "E50" "'  +3.34%  "
Set-CellText $ws "D51" "'144.43"
Set-CellText $ws "E51" "'  +2.59%  "
